# Generate Report for Handoff
# Adds a new tracked file (f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md) to the
# localization-status workbook: one new row on the "Overview" sheet and one
# new row on each of the "zh-cn" / "de-de" language sheets, each backed by
# its ListObject table (so the table ref / autofilter / dimension grow with
# the data), plus the matching hyperlinks.

$wb = $excel.ActiveWorkbook

$commitSha = "db79f4033174aad60540b378c64f25543e341306"
$fileName = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.md"
$pathAndName = "e2e\$fileName"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$fileName"
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $pathAndName) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").NumberFormat = $dateFormat
$wsOverview.Range("G3").Value = "2016-08-25 02:39:33"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $fileName) | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.0e0ee87858c68783cd4d11057245d5d9d0c48721.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("H3").Value = "2016-08-25 02:39:28"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $fileName) | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "f04efb91-24d9-4f68-89a4-3b0ba9fc450c.0e0ee87858c68783cd4d11057245d5d9d0c48721.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("H3").Value = "2016-08-25 02:39:33"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"
